# Auto-generated script applying scheduled-runner price/profit updates
# to the Tonberry_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1482.8334  # H32: was 1532.75
$ws.Cells.Item(32, 10).Value = 1320.7142  # J32: was 1406.2858
$ws.Cells.Item(32, 12).Value = 1320.7142  # L32: was 1406.2858
$ws.Cells.Item(32, 14).Value = -1972.7142  # N32: was -2058.2858

$ws.Cells.Item(62, 8).Value = 1998  # H62: was 2000
$ws.Cells.Item(62, 9).Value = 1995  # I62: was 0
$ws.Cells.Item(62, 11).Value = 1995  # K62: was 0
$ws.Cells.Item(62, 13).Value = -1371  # M62: was None

$ws.Cells.Item(65, 8).Value = 1998  # H65: was 2000
$ws.Cells.Item(65, 9).Value = 1995  # I65: was 0
$ws.Cells.Item(65, 11).Value = 9975  # K65: was 0
$ws.Cells.Item(65, 13).Value = -6855  # M65: was None

$ws.Cells.Item(80, 8).Value = 4680.2  # H80: was 4400.3335
$ws.Cells.Item(80, 10).Value = 1999  # J80: was 2500
$ws.Cells.Item(80, 12).Value = 5997  # L80: was 7500
$ws.Cells.Item(80, 14).Value = -7993  # N80: was -9496

$ws.Cells.Item(83, 8).Value = 4680.2  # H83: was 4400.3335
$ws.Cells.Item(83, 10).Value = 1999  # J83: was 2500
$ws.Cells.Item(83, 12).Value = 17991  # L83: was 22500
$ws.Cells.Item(83, 14).Value = -27975  # N83: was -32484

$ws.Cells.Item(86, 8).Value = 1899.5  # H86: was 1999
$ws.Cells.Item(86, 9).Value = 0  # I86: was 1999
$ws.Cells.Item(86, 10).Value = 1899.5  # J86: was 0
$ws.Cells.Item(86, 11).Value = 0  # K86: was 1999
$ws.Cells.Item(86, 12).Value = 1899.5  # L86: was 0
$ws.Cells.Item(86, 13).ClearContents()  # M86: was -876
$ws.Cells.Item(86, 14).Value = -4145.5  # N86: was None

$ws.Cells.Item(89, 8).Value = 1899.5  # H89: was 1999
$ws.Cells.Item(89, 9).Value = 0  # I89: was 1999
$ws.Cells.Item(89, 10).Value = 1899.5  # J89: was 0
$ws.Cells.Item(89, 11).Value = 0  # K89: was 9995
$ws.Cells.Item(89, 12).Value = 9497.5  # L89: was 0
$ws.Cells.Item(89, 13).ClearContents()  # M89: was -4379
$ws.Cells.Item(89, 14).Value = -20729.5  # N89: was None

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 232308.83  # H2: was 232317.12
$ws.Cells.Item(2, 9).Value = 347882.2  # I2: was 327447.88
$ws.Cells.Item(2, 10).Value = 1162.125  # J2: was 1285.2858
$ws.Cells.Item(2, 11).Value = 347882.2  # K2: was 327447.88
$ws.Cells.Item(2, 12).Value = 1162.125  # L2: was 1285.2858
$ws.Cells.Item(2, 13).Value = -347769.2  # M2: was -327334.88
$ws.Cells.Item(2, 14).Value = -1388.125  # N2: was -1511.2858

$ws.Cells.Item(5, 8).Value = 25000200  # H5: was 50000200
$ws.Cells.Item(5, 9).Value = 267  # I5: was 400
$ws.Cells.Item(5, 11).Value = 267  # K5: was 400
$ws.Cells.Item(5, 13).Value = -155  # M5: was -288

$ws.Cells.Item(32, 8).Value = 3847.9485  # H32: was 4478.7
$ws.Cells.Item(32, 9).Value = 3847.9485  # I32: was 4495.687
$ws.Cells.Item(32, 10).Value = 0  # J32: was 2797
$ws.Cells.Item(32, 11).Value = 3847.9485  # K32: was 4495.687
$ws.Cells.Item(32, 12).Value = 0  # L32: was 2797
$ws.Cells.Item(32, 13).Value = -3560.9485  # M32: was -4208.687
$ws.Cells.Item(32, 14).ClearContents()  # N32: was -3371

$ws.Cells.Item(116, 8).Value = 232308.83  # H116: was 232317.12
$ws.Cells.Item(116, 9).Value = 347882.2  # I116: was 327447.88
$ws.Cells.Item(116, 10).Value = 1162.125  # J116: was 1285.2858
$ws.Cells.Item(116, 11).Value = 347882.2  # K116: was 327447.88
$ws.Cells.Item(116, 12).Value = 1162.125  # L116: was 1285.2858
$ws.Cells.Item(116, 13).Value = -345588.2  # M116: was -325153.88
$ws.Cells.Item(116, 14).Value = -5750.125  # N116: was -5873.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 232308.83  # H3: was 232317.12
$ws.Cells.Item(3, 9).Value = 347882.2  # I3: was 327447.88
$ws.Cells.Item(3, 10).Value = 1162.125  # J3: was 1285.2858
$ws.Cells.Item(3, 11).Value = 347882.2  # K3: was 327447.88
$ws.Cells.Item(3, 12).Value = 1162.125  # L3: was 1285.2858
$ws.Cells.Item(3, 13).Value = -347768.2  # M3: was -327333.88
$ws.Cells.Item(3, 14).Value = -1390.125  # N3: was -1513.2858

$ws.Cells.Item(4, 8).Value = 25000200  # H4: was 50000200
$ws.Cells.Item(4, 9).Value = 267  # I4: was 400
$ws.Cells.Item(4, 11).Value = 267  # K4: was 400
$ws.Cells.Item(4, 13).Value = -152  # M4: was -285

$ws.Cells.Item(82, 8).Value = 21101  # H82: was 13702
$ws.Cells.Item(82, 10).Value = 32400  # J82: was 48000
$ws.Cells.Item(82, 12).Value = 32400  # L82: was 48000
$ws.Cells.Item(82, 14).Value = -33166  # N82: was -48766

$ws.Cells.Item(85, 8).Value = 21101  # H85: was 13702
$ws.Cells.Item(85, 10).Value = 32400  # J85: was 48000
$ws.Cells.Item(85, 12).Value = 32400  # L85: was 48000
$ws.Cells.Item(85, 14).Value = -35052  # N85: was -50652

$ws.Cells.Item(86, 8).Value = 334356.16  # H86: was 445507.66
$ws.Cells.Item(86, 9).Value = 1101  # I86: was 1233.3334
$ws.Cells.Item(86, 10).Value = 572395.5600000001  # J86: was 667644.8
$ws.Cells.Item(86, 11).Value = 1101  # K86: was 1233.3334
$ws.Cells.Item(86, 12).Value = 572395.5600000001  # L86: was 667644.8
$ws.Cells.Item(86, 13).Value = 22  # M86: was -110.3334
$ws.Cells.Item(86, 14).Value = -574641.5600000001  # N86: was -669890.8

$ws.Cells.Item(89, 8).Value = 334356.16  # H89: was 445507.66
$ws.Cells.Item(89, 9).Value = 1101  # I89: was 1233.3334
$ws.Cells.Item(89, 10).Value = 572395.5600000001  # J89: was 667644.8
$ws.Cells.Item(89, 11).Value = 5505  # K89: was 6166.666999999999
$ws.Cells.Item(89, 12).Value = 2861977.8  # L89: was 3338224
$ws.Cells.Item(89, 13).Value = 111  # M89: was -550.6669999999995
$ws.Cells.Item(89, 14).Value = -2873209.8  # N89: was -3349456

$ws.Cells.Item(94, 8).Value = 263.4  # H94: was 292
$ws.Cells.Item(94, 9).Value = 279.25  # I94: was 322.66666
$ws.Cells.Item(94, 11).Value = 279.25  # K94: was 322.66666
$ws.Cells.Item(94, 13).Value = 171.75  # M94: was 128.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1270.2142  # H22: was 1342
$ws.Cells.Item(22, 9).Value = 396.25  # I22: was 430
$ws.Cells.Item(22, 10).Value = 1619.8  # J22: was 1590.7273
$ws.Cells.Item(22, 11).Value = 396.25  # K22: was 430
$ws.Cells.Item(22, 12).Value = 1619.8  # L22: was 1590.7273
$ws.Cells.Item(22, 13).Value = -46.25  # M22: was -80
$ws.Cells.Item(22, 14).Value = -2319.8  # N22: was -2290.7273

$ws.Cells.Item(95, 8).Value = 35000  # H95: was 34812
$ws.Cells.Item(95, 10).Value = 35000  # J95: was 34812
$ws.Cells.Item(95, 12).Value = 35000  # L95: was 34812
$ws.Cells.Item(95, 14).Value = -40492  # N95: was -40304

$ws.Cells.Item(134, 8).Value = 1434.72  # H134: was 1506.4546
$ws.Cells.Item(134, 9).Value = 1293.7142  # I134: was 1357.8889
$ws.Cells.Item(134, 11).Value = 3881.1426  # K134: was 4073.6667
$ws.Cells.Item(134, 13).Value = -1346.1426  # M134: was -1538.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 114.15  # H12: was 114
$ws.Cells.Item(12, 10).Value = 130.53334  # J12: was 130.33333
$ws.Cells.Item(12, 12).Value = 391.60002  # L12: was 390.99999
$ws.Cells.Item(12, 14).Value = -737.6000200000001  # N12: was -736.99999

$ws.Cells.Item(16, 8).Value = 166671680  # H16: was 111114580
$ws.Cells.Item(16, 9).Value = 166671680  # I16: was 111114580
$ws.Cells.Item(16, 11).Value = 500015040  # K16: was 333343740
$ws.Cells.Item(16, 13).Value = -500014867  # M16: was -333343567

$ws.Cells.Item(20, 8).Value = 1874.75  # H20: was 1799.8
$ws.Cells.Item(20, 10).Value = 2749.5  # J20: was 2333
$ws.Cells.Item(20, 12).Value = 8248.5  # L20: was 6999
$ws.Cells.Item(20, 14).Value = -8702.5  # N20: was -7453

$ws.Cells.Item(21, 8).Value = 1725.25  # H21: was 1700.25
$ws.Cells.Item(21, 9).Value = 0  # I21: was 900
$ws.Cells.Item(21, 10).Value = 1725.25  # J21: was 1967
$ws.Cells.Item(21, 11).Value = 0  # K21: was 2700
$ws.Cells.Item(21, 12).Value = 5175.75  # L21: was 5901
$ws.Cells.Item(21, 13).ClearContents()  # M21: was -2527
$ws.Cells.Item(21, 14).Value = -5521.75  # N21: was -6247

$ws.Cells.Item(22, 8).Value = 335600  # H22: was 252112.38
$ws.Cells.Item(22, 9).Value = 501500  # I22: was 334466.34
$ws.Cells.Item(22, 10).Value = 252650  # J22: was 202700
$ws.Cells.Item(22, 11).Value = 1504500  # K22: was 1003399.02
$ws.Cells.Item(22, 12).Value = 757950  # L22: was 608100
$ws.Cells.Item(22, 13).Value = -1504331  # M22: was -1003230.02
$ws.Cells.Item(22, 14).Value = -758288  # N22: was -608438

$ws.Cells.Item(27, 8).Value = 335600  # H27: was 252112.38
$ws.Cells.Item(27, 9).Value = 501500  # I27: was 334466.34
$ws.Cells.Item(27, 10).Value = 252650  # J27: was 202700
$ws.Cells.Item(27, 11).Value = 1504500  # K27: was 1003399.02
$ws.Cells.Item(27, 12).Value = 757950  # L27: was 608100
$ws.Cells.Item(27, 13).Value = -1504398  # M27: was -1003297.02
$ws.Cells.Item(27, 14).Value = -758154  # N27: was -608304

$ws.Cells.Item(33, 8).Value = 199.46153  # H33: was 294.7857
$ws.Cells.Item(33, 9).Value = 138.2  # I33: was 126.833336
$ws.Cells.Item(33, 10).Value = 237.75  # J33: was 420.75
$ws.Cells.Item(33, 11).Value = 829.1999999999999  # K33: was 761.000016
$ws.Cells.Item(33, 12).Value = 1426.5  # L33: was 2524.5
$ws.Cells.Item(33, 13).Value = -546.1999999999999  # M33: was -478.000016
$ws.Cells.Item(33, 14).Value = -1992.5  # N33: was -3090.5

$ws.Cells.Item(80, 8).Value = 2500  # H80: was 2285.2856
$ws.Cells.Item(80, 10).Value = 2500  # J80: was 2285.2856
$ws.Cells.Item(80, 12).Value = 7500  # L80: was 6855.8568
$ws.Cells.Item(80, 14).Value = -9372  # N80: was -8727.856800000001

$ws.Cells.Item(83, 8).Value = 2500  # H83: was 2285.2856
$ws.Cells.Item(83, 10).Value = 2500  # J83: was 2285.2856
$ws.Cells.Item(83, 12).Value = 22500  # L83: was 20567.5704
$ws.Cells.Item(83, 14).Value = -31860  # N83: was -29927.5704

$ws.Cells.Item(113, 8).Value = 1046.2  # H113: was 1065.7587
$ws.Cells.Item(113, 10).Value = 1044.1538  # J113: was 1066.76
$ws.Cells.Item(113, 12).Value = 3132.4614  # L113: was 3200.28
$ws.Cells.Item(113, 14).Value = -7472.4614  # N113: was -7540.28

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1764.4736  # H97: was 1769.6842
$ws.Cells.Item(97, 9).Value = 1687.8572  # I97: was 1694.9286
$ws.Cells.Item(97, 11).Value = 1687.8572  # K97: was 1694.9286
$ws.Cells.Item(97, 13).Value = -1191.8572  # M97: was -1198.9286

$ws.Cells.Item(132, 8).Value = 1103920.5  # H132: was 1103969.1
$ws.Cells.Item(132, 9).Value = 1608341.9  # I132: was 1608412.8
$ws.Cells.Item(132, 11).Value = 4825025.699999999  # K132: was 4825238.4
$ws.Cells.Item(132, 13).Value = -4822495.699999999  # M132: was -4822708.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 5109.077  # H16: was 5113.154
$ws.Cells.Item(16, 9).Value = 5528.1665  # I16: was 6026.4546
$ws.Cells.Item(16, 10).Value = 80  # J16: was 90
$ws.Cells.Item(16, 11).Value = 5528.1665  # K16: was 6026.4546
$ws.Cells.Item(16, 12).Value = 80  # L16: was 90
$ws.Cells.Item(16, 13).Value = -5358.1665  # M16: was -5856.4546
$ws.Cells.Item(16, 14).Value = -420  # N16: was -430

$ws.Cells.Item(82, 8).Value = 2134.4443  # H82: was 2361.3333
$ws.Cells.Item(82, 9).Value = 1719.2307  # I82: was 1935
$ws.Cells.Item(82, 11).Value = 1719.2307  # K82: was 1935
$ws.Cells.Item(82, 13).Value = -1358.2307  # M82: was -1574

$ws.Cells.Item(85, 8).Value = 2134.4443  # H85: was 2361.3333
$ws.Cells.Item(85, 9).Value = 1719.2307  # I85: was 1935
$ws.Cells.Item(85, 11).Value = 1719.2307  # K85: was 1935
$ws.Cells.Item(85, 13).Value = -471.2307000000001  # M85: was -687

$ws.Cells.Item(132, 8).Value = 4244.48  # H132: was 4138.846
$ws.Cells.Item(132, 9).Value = 2399.5  # I132: was 2339.4
$ws.Cells.Item(132, 11).Value = 7198.5  # K132: was 7018.200000000001
$ws.Cells.Item(132, 13).Value = -4668.5  # M132: was -4488.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 33240  # H16: was 45000
$ws.Cells.Item(16, 10).Value = 33240  # J16: was 45000
$ws.Cells.Item(16, 12).Value = 33240  # L16: was 45000
$ws.Cells.Item(16, 14).Value = -33824  # N16: was -45584

$ws.Cells.Item(75, 8).Value = 47110  # H75: was 0
$ws.Cells.Item(75, 10).Value = 47110  # J75: was 0
$ws.Cells.Item(75, 12).Value = 47110  # L75: was 0
$ws.Cells.Item(75, 14).Value = -48982  # N75: was None

$ws.Cells.Item(78, 8).Value = 47110  # H78: was 0
$ws.Cells.Item(78, 10).Value = 47110  # J78: was 0
$ws.Cells.Item(78, 12).Value = 141330  # L78: was 0
$ws.Cells.Item(78, 14).Value = -150690  # N78: was None

$ws.Cells.Item(128, 8).Value = 32666.666  # H128: was 32105.264
$ws.Cells.Item(128, 10).Value = 32666.666  # J128: was 32105.264
$ws.Cells.Item(128, 12).Value = 32666.666  # L128: was 32105.264
$ws.Cells.Item(128, 14).Value = -42626.666  # N128: was -42065.264

$ws.Cells.Item(131, 8).Value = 70000  # H131: was 45000
$ws.Cells.Item(131, 10).Value = 70000  # J131: was 45000
$ws.Cells.Item(131, 12).Value = 70000  # L131: was 45000
$ws.Cells.Item(131, 14).Value = -80080  # N131: was -55080

$ws.Cells.Item(135, 8).Value = 106931.555  # H135: was 105618.4
$ws.Cells.Item(135, 10).Value = 106931.555  # J135: was 105618.4
$ws.Cells.Item(135, 12).Value = 106931.555  # L135: was 105618.4
$ws.Cells.Item(135, 14).Value = -117071.555  # N135: was -115758.4

$ws.Cells.Item(136, 8).Value = 17362772  # H136: was 15433627
$ws.Cells.Item(136, 9).Value = 26456228  # I136: was 22223306
$ws.Cells.Item(136, 11).Value = 79368684  # K136: was 66669918
$ws.Cells.Item(136, 13).Value = -79366134  # M136: was -66667368

$ws.Cells.Item(138, 8).Value = 0  # H138: was 52495
$ws.Cells.Item(138, 10).Value = 0  # J138: was 52495
$ws.Cells.Item(138, 12).Value = 0  # L138: was 52495
$ws.Cells.Item(138, 14).ClearContents()  # N138: was -62775

$ws.Cells.Item(139, 8).Value = 69436.27  # H139: was 69292.234
$ws.Cells.Item(139, 10).Value = 69436.27  # J139: was 69292.234
$ws.Cells.Item(139, 12).Value = 69436.27  # L139: was 69292.234
$ws.Cells.Item(139, 14).Value = -79716.27  # N139: was -79572.234

$ws.Cells.Item(141, 8).Value = 64915.875  # H141: was 65450.934
$ws.Cells.Item(141, 10).Value = 65376.934  # J141: was 65983.14
$ws.Cells.Item(141, 12).Value = 65376.934  # L141: was 65983.14
$ws.Cells.Item(141, 14).Value = -75736.93400000001  # N141: was -76343.14
